$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 values (was row 3: phenol)
$ws.Range("A2").Value = "phenol"
$ws.Range("B2").Value = "phenol"
$ws.Range("C2").Value = "C6H6O"
$ws.Range("D2").Value = "C1=CC=C(C=C1)O"
$ws.Range("E2").Value = 94.11
$ws.Range("F2").Value = 1.5
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0.765763468281798
$ws.Range("K2").Value = 0.06426522154925088
$ws.Range("L2").Value = 0.1700031877590054
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.8193178195728402
$ws.Range("R2").Value = 0.1807140580172139

# New row 3 values (was row 4: dodecane)
$ws.Range("A3").Value = "dodecane"
$ws.Range("B3").Value = "dodecane"
$ws.Range("C3").Value = "C12H26"
$ws.Range("D3").Value = "CCCCCCCCCCCC"
$ws.Range("E3").Value = 170.33
$ws.Range("F3").Value = 6.1
$ws.Range("G3").Value = 12
$ws.Range("H3").Value = 26
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0.846192684788352
$ws.Range("K3").Value = 0.1538660247754359
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 12
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 1.000058709563788
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0

# New row 4 values (was row 2: naphthalene)
$ws.Range("A4").Value = "naphthalene"
$ws.Range("B4").Value = "naphthalene"
$ws.Range("C4").Value = "C10H8"
$ws.Range("D4").Value = "C1=CC=C2C=CC=CC2=C1"
$ws.Range("E4").Value = 128.17
$ws.Range("F4").Value = 3.3
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.9371147694468284
$ws.Range("K4").Value = 0.06291643910431459
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 10
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1.000031208551143
$ws.Range("R4").Value = 0
